$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.463.18'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '3.554.96'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '606.25'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '144.95'
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('D7').Value = '3.553.76'
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.496'
$ws.Range('E9').Value = '  +3.33%  '
$ws.Range('E10').Value = '  -0.80%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.97'
$ws.Range('E11').Value = '  -2.16%  '
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('D13').Value = '4.161.64'
$ws.Range('E13').Value = '  +0.62%  '
$ws.Range('E14').Value = '  -0.36%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '30.00'
$ws.Range('E15').Value = '  -0.91%  '
$ws.Range('D16').Value = '3.540.31'
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('D17').Value = '66.547.35'
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.115'
$ws.Range('E18').Value = '  +0.18%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.50'
$ws.Range('E19').Value = '  +5.11%  '
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.91'
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '430.94'
$ws.Range('E22').Value = '  +0.99%  '
$ws.Range('E23').Value = '  +2.07%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '79.41'
$ws.Range('E24').Value = '  +0.61%  '
$ws.Range('D25').Value = '3.698.44'
$ws.Range('E25').Value = '  +0.63%  '
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000119'
$ws.Range('E27').Value = '  -0.67%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.51'
$ws.Range('E28').Value = '  +0.93%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.99'
$ws.Range('E29').Value = '  -1.81%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '9.14'
$ws.Range('E30').Value = '  -1.16%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('E32').Value = '  -2.17%  '
$ws.Range('D33').Value = '3.552.64'
$ws.Range('E33').Value = '  +0.72%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '25.31'
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.153'
$ws.Range('E35').Value = '  -4.07%  '
$ws.Range('B36').Value = 'USDe'
$ws.Range('C36').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '7.81'
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('E38').Value = '  -2.01%  '
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '172.90'
$ws.Range('E40').Value = '  -0.35%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0846'
$ws.Range('E41').Value = '  -1.36%  '
$ws.Range('E42').Value = '  -1.44%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.886'
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.92'
$ws.Range('E44').Value = '  +0.85%  '
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.50'
$ws.Range('E46').Value = '  +3.91%  '
$ws.Range('E47').Value = '  -2.80%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '25.04'
$ws.Range('E48').Value = '  -3.87%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.16'
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '23.48'
$ws.Range('E50').Value = '  +4.13%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.941'
$ws.Range('E51').Value = '  -0.35%  '
